$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-12 Monday" "2024-02-13 Tuesday"

Replace-Text "93÷3=" "93÷4="
Replace-Text "65÷3=" "86÷5="
Replace-Text "81÷8=" "95÷3="
Replace-Text "56÷4=" "78÷4="
Replace-Text "12÷3=" "51÷9="

Replace-Text "98÷9=" "97÷6="
Replace-Text "10÷8=" "18÷8="
Replace-Text "15÷3=" "85÷6="
Replace-Text "16÷9=" "53÷8="
Replace-Text "75÷4=" "64÷9="

Replace-Text "14÷2=" "57÷6="
Replace-Text "43÷3=" "80÷8="
Replace-Text "69÷4=" "29÷2="
Replace-Text "96÷5=" "97÷8="
Replace-Text "28÷8=" "37÷8="

Replace-Text "87÷2=" "90÷4="
Replace-Text "89÷9=" "50÷6="
Replace-Text "80÷5=" "33÷3="
Replace-Text "42÷9=" "76÷4="
Replace-Text "85÷9=" "53÷3="

Replace-Text "78÷2=" "80÷2="
Replace-Text "56÷8=" "55÷7="
Replace-Text "59÷6=" "95÷5="
Replace-Text "70÷3=" "33÷4="
Replace-Text "87÷6=" "98÷5="
